$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3: date-looking text must stay plain text (like A2), not auto-converted to a date serial.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-11-02"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = 0.7616438356164383

$ws.Range("C3").Value = "{'c_w': 0.25447054029049226, 'c_s': 1.8568662619089953, 'den_mid': 0.46439117657448753, 'den_high': 0.6446772564277745, 'low_thres': 0.6399957717850465, 'high_thres': 0.6979598879442052, 'p_base_h': 0.20097534714976714, 'c_r': 1.1609414221653496, 'c_t': 1.082500226109699, 't_ref_c': 33.17811951296828, 't_scale': 5.030357500017056, 'rain_win': 3, 'r_ref_mm': 6.6653330727891, 'yama_step': 62, 'c_y': 0.03582999029214878, 'yama_cap': 2.386033787257072, 'burn_steps': 6, 'ember_chance': 0.04340509960691442, 'c_h': 0.024262980583097227, 'h_ref': 79.14251040914458, 'h_scale': 17.568516859786204, 'c_sm': 0.7011653727056256}"

$ws.Range("D3").Value = "Đã thêm độ ẩm không khí (c_h, h_ref, h_scale) và độ ẩm đất (c_sm). So sánh F1 theo trạng thái cuối sau đúng n_t, không dừng sớm."
